$d = $word.ActiveDocument

# Trim the trailing clause from the SmartCash mining paragraph so it now
# ends right after "quite some time." instead of continuing on about
# Smartcash reaching a considerable market cap.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(", until Smartcash reaches a considerable market cap.", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)
